# Insert a new weekly record at row 52 (pushing the existing rows 52-110
# down to 53-111) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 52:110 down to 53:111, carrying formatting (incl. the date
# style on column D) along with them.
$ws.Rows(52).Insert()

# Populate the newly inserted row 52 with the new record.
$ws.Range("A52").Value = 7
$ws.Range("B52").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C52").Value = "Ñuble"
$ws.Range("D52").Value = 44539
$ws.Range("E52").Value = 16
$ws.Range("F52").Value = 100112024
$ws.Range("G52").Value = "Choclo"
$ws.Range("H52").Value = "Choclero"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 800
$ws.Range("K52").Value = 450
$ws.Range("L52").Value = 500
$ws.Range("M52").Value = 475
$ws.Range("N52").Value = "$/unidad"
$ws.Range("O52").Value = "Región Metropolitana"
$ws.Range("P52").Value = 475
$ws.Range("Q52").Value = 1
$ws.Range("R52").Value = "Hortaliza"
